# JSF - 5 - Conversores y validadores: mejoras pdfs y codigo de conversores
#
# 1) Date-stamp placeholders ("datetimeFigureOut" fields) in the three
#    slide masters and the layouts that own a live copy of the field:
#    18/10/2022 -> 31/10/2022
# 2) Slide 3: "convertidores" -> "conversores" wording fix.
# 3) Slide 5: resize/reposition the JSF snippet textbox and move the
#    converterMessage attribute line up into the first tag so the XML
#    sample reads as a single <h:inputText ...> element.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Fix the cached date field text on every master / layout that has it
# ---------------------------------------------------------------------
for ($di = 1; $di -le $p.Designs.Count; $di++) {
    $design = $p.Designs.Item($di)
    $master = $design.SlideMaster

    for ($si = 1; $si -le $master.Shapes.Count; $si++) {
        $sh = $master.Shapes.Item($si)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "18/10/2022") {
                $sh.TextFrame.TextRange.Text = "31/10/2022"
            }
        }
    }

    for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
        $layout = $master.CustomLayouts.Item($li)
        for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
            $sh = $layout.Shapes.Item($si)
            if ($sh.HasTextFrame) {
                $t = $sh.TextFrame.TextRange.Text
                if ($t -eq "18/10/2022") {
                    $sh.TextFrame.TextRange.Text = "31/10/2022"
                }
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 3 - "convertidores" -> "conversores"
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$contentShape = $slide3.Shapes.Item(2)
$contentShape.TextFrame.TextRange.Paragraphs(1).Text = "Se pueden aplicar conversores a las siguientes etiquetas"

# ---------------------------------------------------------------------
# 3) Slide 5 - resize textbox + move converterMessage line into tag 1
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$box = $slide5.Shapes.Item(4)

$box.Left = 25.474645669291338
$box.Top = 207.63086614173227
$box.Width = 658.5252755905511
$box.Height = 66.90913385826772

$tr = $box.TextFrame.TextRange

# Paragraph 1 currently ends right before its paragraph mark with the
# run '}" '. Insert the new attribute text immediately after it - the
# 6 inserts below need to happen back-to-front against the (fixed)
# paragraph-mark position so the resulting text lands in forward order.
$mark = $tr.Characters(53, 1)
$mark.InsertBefore(" favor!`"> `t") | Out-Null
$mark.InsertBefore("por") | Out-Null
$mark.InsertBefore(" ") | Out-Null
$mark.InsertBefore("Entero") | Out-Null
$mark.InsertBefore("=" + [char]0x201C) | Out-Null
$mark.InsertBefore("converterMessage") | Out-Null

# Apply matching run formatting (Gill Sans 17pt black) to the whole
# freshly-inserted span.
$newSpan = $tr.Characters(53, 39)
$newSpan.Font.Name = "Gill Sans"
$newSpan.Font.Size = 17
$newSpan.Font.Color.RGB = 0

# Paragraph 2 used to start with a second tab + the converterMessage
# attribute text that has just been relocated above; trim it back down
# to the trailing tab + "<" that starts the <f:convertNumber .../> tag.
$tr2 = $box.TextFrame.TextRange
$tr2.Characters(92, 43).Text = ""
